$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.331.53"
$ws.Range("E2").Value = "  +0.54%  "
$ws.Range("D3").Value = "2.550.04"
$ws.Range("E3").Value = "  -2.15%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "590.82"
$ws.Range("E5").Value = "  +0.17%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "174.44"
$ws.Range("E6").Value = "  +5.72%  "
$ws.Range("E7").Value = "  +0.07%  "
$ws.Range("E8").Value = "  -0.30%  "
$ws.Range("D9").Value = "2.549.62"
$ws.Range("E9").Value = "  -2.16%  "
$ws.Range("E10").Value = "  +1.69%  "
$ws.Range("E11").Value = "  +1.24%  "
$ws.Range("E12").Value = "  -2.36%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.16"
$ws.Range("E13").Value = "  -0.37%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "27.10"
$ws.Range("E14").Value = "  -0.44%  "
$ws.Range("D15").Value = "3.011.56"
$ws.Range("E15").Value = "  -2.24%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0000178"
$ws.Range("E16").Value = "  -0.13%  "
$ws.Range("D17").Value = "67.269.73"
$ws.Range("E17").Value = "  +0.33%  "
$ws.Range("D18").Value = "2.566.45"
$ws.Range("E18").Value = "  -1.91%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "8.07"
$ws.Range("E19").Value = "  +3.68%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.41"
$ws.Range("E20").Value = "  -2.92%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "356.16"
$ws.Range("E21").Value = "  +0.82%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.23"
$ws.Range("E22").Value = "  -0.73%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.68"
$ws.Range("E23").Value = "  +1.50%  "
$ws.Range("E24").Value = "  +3.39%  "
$ws.Range("E25").Value = "  +0.08%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "70.14"
$ws.Range("E26").Value = "  +1.83%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.97"
$ws.Range("E27").Value = "  -5.19%  "
$ws.Range("E28").Value = "  -2.29%  "
$ws.Range("E29").Value = "  +0.22%  "
$ws.Range("E30").Value = "  +0.94%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "536.38"
$ws.Range("E31").Value = "  -0.43%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "8.28"
$ws.Range("E32").Value = "  +5.45%  "
$ws.Range("E33").Value = "  +1.23%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.86"
$ws.Range("E34").Value = "  -0.35%  "
$ws.Range("E35").Value = "  -0.85%  "
$ws.Range("E36").Value = "  +0.09%  "
$ws.Range("E37").Value = "  -0.06%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "157.67"
$ws.Range("E38").Value = "  +0.43%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "18.82"
$ws.Range("E39").Value = "  -0.23%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "18.43"
$ws.Range("E40").Value = "  +1.16%  "
$ws.Range("E41").Value = "  -1.60%  "
$ws.Range("E42").Value = "  +0.66%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.20"
$ws.Range("E43").Value = "  +1.66%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.58"
$ws.Range("E44").Value = "  +7.74%  "
$ws.Range("E45").Value = "  +0.03%  "
$ws.Range("E46").Value = "  -0.64%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "152.22"
$ws.Range("E47").Value = "  +0.97%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.565"
$ws.Range("E48").Value = "  -1.28%  "
$ws.Range("E49").Value = "  -5.81%  "
$ws.Range("E50").Value = "  -0.84%  "
$ws.Range("E51").Value = "  +1.85%  "
